$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 278. This shifts the existing rows
# 278-338 down to 280-340, preserving all of their data and formatting.
$ws.Rows.Item(278).Insert()
$ws.Rows.Item(278).Insert()

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

# New data for row 278 (Primera quality entry)
$row278 = @(5, "Macroferia Regional de Talca", "Maule", 44785, 7, 100112006, "Repollo", "Crespo record", "Primera", 3000, 1200, 1200, 1200, "`$/unidad", "Región del Maule", 1200, 1, "Hortaliza")

# New data for row 279 (Segunda quality entry)
$row279 = @(5, "Macroferia Regional de Talca", "Maule", 44785, 7, 100112006, "Repollo", "Crespo record", "Segunda", 3000, 900, 900, 900, "`$/unidad", "Región del Maule", 900, 1, "Hortaliza")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "278").Value = $row278[$i]
    $ws.Range($cols[$i] + "279").Value = $row279[$i]
}
